$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the short name for the loan product
$ws.Range("B3").Value = 919

# Fix the swapped "decimal places" / "currency" values:
#  - decimal places should be the number 2
#  - currency should be the text "US Dollar"
$ws.Range("B6").Value = 2
$ws.Range("B7").Value = "US Dollar"

# Leave the selection on the currency row that was just edited
$ws.Range("A7:B7").Select()
